{"js": "// \"Draw the game board\" -> \"Draws the game board\"\n//\n// The canonical edit splits the original single run into three runs\n// (\"Draw\" / \"s\" / \" the game board\") with the document's \"_GoBack\"\n// bookmark (which Word maintains at the location of the most recent\n// edit) now sitting right after the inserted \"s\". The \"_GoBack\"\n// bookmark that used to sit at the very end of the document (after\n// \"Performs computer player move\") is removed, since Word only ever\n// keeps a single \"_GoBack\" bookmark.\n\nconst body = context.document.body;\n\n// 1) Remove the stale \"_GoBack\" bookmark from the end of the document.\n//    At this point it is the only bookmark with that name, so this\n//    unambiguously clears the old location.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// 2) Locate \"Draw the game board\" and replace it in-place with the same\n//    text split across three runs, with a new \"_GoBack\" bookmark\n//    inserted between the \"s\" and the remaining \" the game board\" text.\nconst matches = body.search(\"Draw the game board\", { matchCase: true });\nmatches.load(\"text\");\nawait context.sync();\n\nif (matches.items.length === 0) {\n  throw new Error('Could not find \"Draw the game board\" in the document.');\n}\n\nconst target = matches.items[0];\n\nconst rPr =\n  '<w:rPr><w:rFonts w:ascii=\"Times New Roman\" w:hAnsi=\"Times New Roman\" w:cs=\"Times New Roman\"/></w:rPr>';\n\nconst ooxml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:r>${rPr}<w:t>Draw</w:t></w:r>\n            <w:r>${rPr}<w:t>s</w:t></w:r>\n            <w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>\n            <w:bookmarkEnd w:id=\"0\"/>\n            <w:r>${rPr}<w:t xml:space=\"preserve\"> the game board</w:t></w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\n\ntarget.insertOoxml(ooxml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# \"Draw the game board\" -> \"Draws the game board\"\n#\n# The canonical edit splits the original single run into three runs\n# (\"Draw\" / \"s\" / \" the game board\") with the document's \"_GoBack\"\n# bookmark (which Word maintains at the location of the most recent\n# edit) now sitting right after the inserted \"s\". The \"_GoBack\"\n# bookmark that used to sit at the very end of the document (after\n# \"Performs computer player move\") is removed, since Word only ever\n# keeps a single \"_GoBack\" bookmark.\n\n$d = $word.ActiveDocument\n\n# 1) Remove the stale \"_GoBack\" bookmark from the end of the document.\n#    At this point it is the only bookmark with that name, so this\n#    unambiguously clears the old location.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# 2) Locate \"Draw the game board\" and replace it in-place with the same\n#    text split across three runs, with a new \"_GoBack\" bookmark\n#    inserted between the \"s\" and the remaining \" the game board\" text.\n$find = $d.Content\n$find.Find.Execute(\"Draw the game board\")\n\n$startPos = $find.Start\n$endPos = $find.End\n\n$target = $d.Range($startPos, $endPos)\n$target.Text = \"\"\n\n$insertionPoint = $d.Range($startPos, $startPos)\n\n$pAttrs = 'w:rsidR=\"00446F72\" w:rsidRDefault=\"00446F72\" w:rsidP=\"00446F72\"'\n$pPr = '<w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"1\"/><w:numId w:val=\"1\"/></w:numPr><w:rPr><w:rFonts w:ascii=\"Times New Roman\" w:hAnsi=\"Times New Roman\" w:cs=\"Times New Roman\"/></w:rPr></w:pPr>'\n$rPr = '<w:rPr><w:rFonts w:ascii=\"Times New Roman\" w:hAnsi=\"Times New Roman\" w:cs=\"Times New Roman\"/></w:rPr>'\n\n$ooxml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n      '<pkg:xmlData>' +\n        '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n          '<w:body>' +\n            '<w:p ' + $pAttrs + '>' +\n              $pPr +\n              '<w:r>' + $rPr + '<w:t>Draw</w:t></w:r>' +\n              '<w:r>' + $rPr + '<w:t>s</w:t></w:r>' +\n              '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>' +\n              '<w:bookmarkEnd w:id=\"0\"/>' +\n              '<w:r>' + $rPr + '<w:t xml:space=\"preserve\"> the game board</w:t></w:r>' +\n            '</w:p>' +\n          '</w:body>' +\n        '</w:document>' +\n      '</pkg:xmlData>' +\n    '</pkg:part>' +\n  '</pkg:package>'\n\n$insertionPoint.InsertXML($ooxml)\n"}
